$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.189.80"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.829.45"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.83"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("E6").Value = "  -3.63%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07121"
$ws.Range("E8").Value = "  -4.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2832"
$ws.Range("E9").Value = "  -2.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.10"
$ws.Range("E10").Value = "  -2.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07653"
$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("D12").Value = "1.827.09"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.794"
$ws.Range("E13").Value = "  -3.77%  "

$ws.Range("E14").Value = "  -5.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001002"
$ws.Range("E15").Value = "  -2.00%  "

$ws.Range("D16").Value = "2.082.54"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.61"
$ws.Range("E17").Value = "  -2.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.025"
$ws.Range("E18").Value = "  -3.57%  "

$ws.Range("D19").Value = "29.174.35"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.27"
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.79"
$ws.Range("E21").Value = "  -4.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.044"
$ws.Range("E23").Value = "  -5.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.052"
$ws.Range("E26").Value = "  -4.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1290"
$ws.Range("E27").Value = "  -4.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.70"
$ws.Range("E28").Value = "  -4.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06891"
$ws.Range("E29").Value = "  +5.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.461"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.456"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.839"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.810"
$ws.Range("E33").Value = "  -6.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.140"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.722"
$ws.Range("E35").Value = "  -6.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6623"
$ws.Range("E36").Value = "  -4.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.533"

$ws.Range("D38").Value = "1.235.73"
$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.754"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01769"
$ws.Range("E40").Value = "  -4.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.587"
$ws.Range("E41").Value = "  -2.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9319"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "1.990.63"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.32"
$ws.Range("E45").Value = "  -0.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.53"
$ws.Range("E46").Value = "  -3.20%  "

$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.641"
$ws.Range("E48").Value = "  -4.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.579"
$ws.Range("E49").Value = "  -6.73%  "

# Row 50 and 51: EnergySwap/Cronos order swapped with updated values
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.536"
$ws.Range("E50").Value = "  -4.83%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05592"
$ws.Range("E51").Value = "  -1.36%  "

